# Apply the commit's data change to the Colos worksheet:
#   - Insert a new row for "Bordeaux, France" (BOD) before the existing
#     "Americana, Brazil" (QWJ) row, i.e. at sheet row 166.
#   - Insert a new row for "San Antonio, United States" (SAT) before the
#     existing "Adelaide, SA, Australia" (ADL) row. After the first
#     insertion that row has moved down to sheet row 291.
# Both insertions shift everything below them down by one row, growing the
# sheet from 301 data+header rows (A1:G301) to 303 (A1:G303).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-ColoRow {
    param($Row, $Colo, $Name, $Lat, $Lon, $Cca2, $Region, $City)

    # Insert a blank row, shifting this row (and everything below) down.
    $ws.Rows.Item($Row).Insert()

    $ws.Range("A$Row").Value = $Colo
    $ws.Range("B$Row").Value = $Name
    $ws.Range("C$Row").Value = $Lat
    $ws.Range("D$Row").Value = $Lon
    $ws.Range("E$Row").Value = $Cca2
    $ws.Range("F$Row").Value = $Region
    $ws.Range("G$Row").Value = $City

    # Column A throughout the table uses the bold / bordered / centered-top
    # style (same as the header row) -- match it on the new row too.
    $colA = $ws.Range("A$Row")
    $colA.Font.Bold = $true
    $colA.Borders.LineStyle = 1
    $colA.HorizontalAlignment = -4108
    $colA.VerticalAlignment = -4160
}

# 1) New row 166: Bordeaux, France (pushes the old row 166 "QWJ" down to 167)
Set-ColoRow 166 "BOD" "Bordeaux, France" 44.82946 -0.58355 "FR" "Europe" "Bordeaux"

# 2) New row 291: San Antonio, United States (the old "ADL" row, originally
#    at 290, is now at 291 after the first insertion above).
Set-ColoRow 291 "SAT" "San Antonio, United States" 29.429461 -98.487061 "US" "North America" "San Antonio"
